# Weekly Reports and Timesheet Tracking -- Summer Week 2 update
# - Mark Summer term "week 1" weekly-report items (formatting/duplicates/mixed-up/
#   desc-fine checks) as reviewed, adding a note about dates/detail still needed.
# - Mark Summer term "week 2" weekly-report row fully checked off.
# - Mark the matching Timesheets rows (week 1 + week 2) as checked off too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$checkMark = [char]0x2714
$green = 65280          # RGB(0,255,0) as a BGR COM long
$xlLeft = -4131
$xlCenter = -4108

function Set-Check($addr) {
    $c = $ws.Range($addr)
    $c.Value = $checkMark
    $c.Font.Name = "Zapf Dingbats"
    $c.Interior.Color = $green
}

function Set-GreenNote($addr, $text) {
    $c = $ws.Range($addr)
    $c.Value = $text
    $c.Interior.Color = $green
}

# --- Weekly Reports Checklist: Summer term, week 1 (row 14) ---
Set-Check "G14"
Set-Check "H14"
Set-Check "I14"
Set-Check "L14"
Set-GreenNote "M14" "detail?"

# New note cell for "dates?" next to the other reviewer comments, merged
# down through rows 12-14 like the existing K7:K14 / L8:L13 notes.
$j12 = $ws.Range("J12")
$j12.Value = "dates?"
$j12.Interior.Color = $green
$j12.HorizontalAlignment = $xlLeft
$j12.VerticalAlignment = $xlCenter
$ws.Range("J12:J14").Merge() | Out-Null

# --- Weekly Reports Checklist: Summer term, week 2 (row 15) -- fully checked ---
Set-Check "D15"
Set-Check "E15"
Set-Check "F15"
Set-Check "G15"
Set-Check "H15"
Set-Check "I15"
Set-Check "J15"
Set-Check "K15"
Set-Check "L15"
Set-Check "M15"

# --- Timesheets Checklist: Summer term, week 1 (row 43) ---
Set-Check "L43"

# --- Timesheets Checklist: Summer term, week 2 (row 44) -- fully checked ---
Set-Check "D44"
Set-Check "E44"
Set-Check "F44"
Set-Check "G44"
Set-Check "H44"
Set-Check "I44"
Set-Check "J44"
Set-Check "K44"
Set-Check "L44"
Set-Check "M44"

# Leave the cursor where the author left it when they saved the file.
$ws.Range("N20").Select() | Out-Null
